# Generate Report for Handoff
# Update localization-status workbook: mark 536304ba-ef74-4d40-a73c-b78d6f8edf69 and
# eff68793-b865-4576-9974-4f7fd12ba60c as "Ready for handoff" with a fresh handoff
# timestamp and a stale-handback-file error message.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

$msg536304ba = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b19dabde65a577d8bd4dc2abe357d1400500a1f/e2e/536304ba-ef74-4d40-a73c-b78d6f8edf69.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2916b70b96fbfbab6a8b03a205d4dee5ec693a8c/e2e/536304ba-ef74-4d40-a73c-b78d6f8edf69.md."
$msgEff68793 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b19dabde65a577d8bd4dc2abe357d1400500a1f/e2e/eff68793-b865-4576-9974-4f7fd12ba60c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2916b70b96fbfbab6a8b03a205d4dee5ec693a8c/e2e/eff68793-b865-4576-9974-4f7fd12ba60c.md."

# ----- Overview sheet -----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E4").Value = $newStatus
$ov.Range("F4").Value = $newStatus
$ov.Range("G4").Value = "2016-08-29 06:26:57"

$ov.Range("E5").Value = $newStatus
$ov.Range("F5").Value = $newStatus
$ov.Range("G5").Value = "2016-08-29 06:26:57"

# ----- zh-cn sheet -----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C4").Value = $newStatus
$zh.Range("H4").Value = "2016-08-29 06:26:53"
$zh.Range("P4").Value = $msg536304ba

$zh.Range("C5").Value = $newStatus
$zh.Range("H5").Value = "2016-08-29 06:26:53"
$zh.Range("P5").Value = $msgEff68793

# Widen the new Error Detail column to match the other wide (40-char) columns.
$zh.Columns.Item(16).ColumnWidth = $zh.Columns.Item(1).ColumnWidth

# ----- de-de sheet -----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C4").Value = $newStatus
$de.Range("H4").Value = "2016-08-29 06:26:57"
$de.Range("P4").Value = $msg536304ba

$de.Range("C5").Value = $newStatus
$de.Range("H5").Value = "2016-08-29 06:26:57"
$de.Range("P5").Value = $msgEff68793

$de.Columns.Item(16).ColumnWidth = $de.Columns.Item(1).ColumnWidth
